# ACO routing relevant input
# Adds a "Sheet2" worksheet after "Sheet1" holding the Node1/Node2/Distance
# table used as the routing-distance input for the ACO example, and leaves
# Sheet1's selection parked at the header row (A1:C1).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- select the header row on Sheet1 before focus moves to the new sheet ---
$ws1.Range("A1:C1").Select()

# --- create Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "Sheet2"

# --- header row ---
$ws2.Range("A1").Value = "Node1"
$ws2.Range("B1").Value = "Node2"
$ws2.Range("C1").Value = "Distance"

# --- edge list (Node1, Node2, Distance) ---
$data = @(
    @(0, 1, 5),
    @(1, 2, 11),
    @(2, 3, 12),
    @(3, 4, 12),
    @(4, 5, 5),
    @(5, 6, 14),
    @(6, 7, 14),
    @(7, 8, 20),
    @(8, 9, 5),
    @(9, 10, 12),
    @(10, 11, 32),
    @(11, 12, 11),
    @(12, 13, 5),
    @(13, 14, 20),
    @(14, 15, 20),
    @(15, 16, 14),
    @(16, 0, 14),
    @(17, 18, 11),
    @(18, 19, 12),
    @(19, 20, 12),
    @(20, 6, 5),
    @(21, 22, 11),
    @(22, 23, 6),
    @(23, 24, 18),
    @(24, 7, 5),
    @(15, 21, 5),
    @(16, 17, 5)
)

$row = 2
foreach ($edge in $data) {
    $ws2.Cells.Item($row, 1).Value = $edge[0]
    $ws2.Cells.Item($row, 2).Value = $edge[1]
    $ws2.Cells.Item($row, 3).Value = $edge[2]
    $row = $row + 1
}

# --- total distance ---
$ws2.Range("C29").Formula = "=SUM(C2:C28)"

# --- leave Sheet2 active with C29 selected, matching the authored file ---
$ws2.Range("C29").Select()
